$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.046598473704398
$ws.Range("D2").Value = 1.054029843438831
$ws.Range("E2").Value = 1.060110850164385
$ws.Range("F2").Value = 1.067292445372106
$ws.Range("I2").Value = 1.048906877990022
$ws.Range("J2").Value = 1.051652013972134
$ws.Range("K2").Value = 1.056774179281239
$ws.Range("L2").Value = 1.062838520788511
$ws.Range("M2").Value = 1.070000699196018
$ws.Range("N2").Value = 1.021076528456047

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.047422920923547
$ws.Range("D3").Value = 1.054674070085734
$ws.Range("E3").Value = 1.060887515031045
$ws.Range("F3").Value = 1.068073110144085
$ws.Range("I3").Value = 1.049130991623463
$ws.Range("J3").Value = 1.052125210303717
$ws.Range("K3").Value = 1.057232000437864
$ws.Range("L3").Value = 1.063429642249805
$ws.Range("M3").Value = 1.070597208603252
$ws.Range("N3").Value = 1.021234844829427

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.047957084737315
$ws.Range("D4").Value = 1.055091494266831
$ws.Range("E4").Value = 1.06139106217238
$ws.Range("F4").Value = 1.068579197512919
$ws.Range("I4").Value = 1.04927506262898
$ws.Range("J4").Value = 1.052431400804257
$ws.Range("K4").Value = 1.057528118318061
$ws.Range("L4").Value = 1.063812467002203
$ws.Range("M4").Value = 1.070983471762295
$ws.Range("N4").Value = 1.021337254019597

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.048181810937662
$ws.Range("D5").Value = 1.055267113148977
$ws.Range("E5").Value = 1.061602989131917
$ws.Range("F5").Value = 1.068792180854712
$ws.Range("I5").Value = 1.049335402859264
$ws.Range("J5").Value = 1.052560122158431
$ws.Range("K5").Value = 1.057652575450274
$ws.Range("L5").Value = 1.063973483832581
$ws.Range("M5").Value = 1.071145922569942
$ws.Range("N5").Value = 1.021380298672809

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.048219553005414
$ws.Range("D6").Value = 1.055296608121588
$ws.Range("E6").Value = 1.06163858636693
$ws.Range("F6").Value = 1.068827954783607
$ws.Range("I6").Value = 1.049345520899992
$ws.Range("J6").Value = 1.05258173493509
$ws.Range("K6").Value = 1.057673470502917
$ws.Range("L6").Value = 1.064000523753205
$ws.Range("M6").Value = 1.071173202584189
$ws.Range("N6").Value = 1.021387525566181

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.047960086898661
$ws.Range("D7").Value = 1.055093840371693
$ws.Range("E7").Value = 1.061393893026988
$ws.Range("F7").Value = 1.068582042528316
$ws.Range("I7").Value = 1.04927586979199
$ws.Range("J7").Value = 1.05243312079141
$ws.Range("K7").Value = 1.057529781443195
$ws.Range("L7").Value = 1.063814618215239
$ws.Range("M7").Value = 1.070985642181568
$ws.Range("N7").Value = 1.021337829217076

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.046876955289803
$ws.Range("D8").Value = 1.054247444445481
$ws.Range("E8").Value = 1.060373120878621
$ws.Range("F8").Value = 1.067556077781445
$ws.Range("I8").Value = 1.048982813493247
$ws.Range("J8").Value = 1.051811931915022
$ws.Range("K8").Value = 1.056928926731713
$ws.Range("L8").Value = 1.063038223830853
$ws.Range("M8").Value = 1.070202233020147
$ws.Range("N8").Value = 1.02113003856676

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.044973712329016
$ws.Range("D9").Value = 1.052760407131366
$ws.Range("E9").Value = 1.05858208360296
$ws.Range("F9").Value = 1.065755524931365
$ws.Range("I9").Value = 1.048459207770325
$ws.Range("J9").Value = 1.050717383013812
$ws.Range("K9").Value = 1.055869262479238
$ws.Range("L9").Value = 1.061672713663891
$ws.Range("M9").Value = 1.068823996785075
$ws.Range("N9").Value = 1.020763660694582

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.043708602276762
$ws.Range("D10").Value = 1.051772131800318
$ws.Range("E10").Value = 1.057393350662956
$ws.Range("F10").Value = 1.064560205080453
$ws.Range("I10").Value = 1.04810534607
$ws.Range("J10").Value = 1.049987807114837
$ws.Range("K10").Value = 1.055162308959807
$ws.Range("L10").Value = 1.060764215172445
$ws.Range("M10").Value = 1.067906771652198
$ws.Range("N10").Value = 1.020519287163563

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.043161699981966
$ws.Range("D11").Value = 1.051344952345233
$ws.Range("E11").Value = 1.056879896556792
$ws.Range("F11").Value = 1.064043841634399
$ws.Range("I11").Value = 1.047950994543677
$ws.Range("J11").Value = 1.049671938748131
$ws.Range("K11").Value = 1.054856086475987
$ws.Range("L11").Value = 1.060371282455079
$ws.Range("M11").Value = 1.067510003277883
$ws.Range("N11").Value = 1.020413448035097

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.042958692938418
$ws.Range("D12").Value = 1.051186393502102
$ws.Range("E12").Value = 1.056689370459286
$ws.Range("F12").Value = 1.063852226334317
$ws.Range("I12").Value = 1.047893493049547
$ws.Range("J12").Value = 1.049554618958092
$ws.Range("K12").Value = 1.054742327139714
$ws.Range("L12").Value = 1.060225399471955
$ws.Range("M12").Value = 1.067362687052438
$ws.Range("N12").Value = 1.020374131594487

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.043002232418097
$ws.Range("D13").Value = 1.051220399694542
$ws.Range("E13").Value = 1.056730230159399
$ws.Range("F13").Value = 1.063893320058016
$ws.Range("I13").Value = 1.047905834924703
$ws.Range("J13").Value = 1.049579784077274
$ws.Range("K13").Value = 1.054766729545005
$ws.Range("L13").Value = 1.060256688685486
$ws.Range("M13").Value = 1.067394284085546
$ws.Range("N13").Value = 1.020382565235682

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.043144916554926
$ws.Range("D14").Value = 1.051331843472443
$ws.Range("E14").Value = 1.056864143648522
$ws.Range("F14").Value = 1.064027998866348
$ws.Range("I14").Value = 1.047946244883853
$ws.Range("J14").Value = 1.049662240891042
$ws.Range("K14").Value = 1.05484668338498
$ws.Range("L14").Value = 1.060359222290329
$ws.Range("M14").Value = 1.067497824814161
$ws.Range("N14").Value = 1.020410198185727

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.04323284714458
$ws.Range("D15").Value = 1.051400522908189
$ws.Range("E15").Value = 1.056946677886553
$ws.Range("F15").Value = 1.064111003511929
$ws.Range("I15").Value = 1.047971120496087
$ws.Range("J15").Value = 1.049713046327544
$ws.Range("K15").Value = 1.054895943682269
$ws.Range("L15").Value = 1.060422405916581
$ws.Range("M15").Value = 1.067561627836624
$ws.Range("N15").Value = 1.020427223365519

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.043744917469105
$ws.Range("D16").Value = 1.051800498264167
$ws.Range("E16").Value = 1.057427453960863
$ws.Range("F16").Value = 1.06459450026105
$ws.Range("I16").Value = 1.048115566200124
$ws.Range("J16").Value = 1.050008771275573
$ws.Range("K16").Value = 1.055182629803507
$ws.Range("L16").Value = 1.060790302510186
$ws.Range("M16").Value = 1.067933112353936
$ws.Range("N16").Value = 1.020526310899377

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.044066367533348
$ws.Range("D17").Value = 1.052051594455837
$ws.Range("E17").Value = 1.057729374918484
$ws.Range("F17").Value = 1.064898112552798
$ws.Range("I17").Value = 1.04820587203454
$ws.Range("J17").Value = 1.050194284139924
$ws.Range("K17").Value = 1.055362432786676
$ws.Range("L17").Value = 1.061021196800568
$ws.Range("M17").Value = 1.068166242038328
$ws.Range("N17").Value = 1.020588459876312

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.044253950508646
$ws.Range("D18").Value = 1.052198126880239
$ws.Range("E18").Value = 1.057905603196394
$ws.Range("F18").Value = 1.065075321838844
$ws.Range("I18").Value = 1.048258437095374
$ws.Range("J18").Value = 1.050302494695053
$ws.Range("K18").Value = 1.055467298412519
$ws.Range("L18").Value = 1.061155917143292
$ws.Range("M18").Value = 1.068302260749074
$ws.Range("N18").Value = 1.020624707995686

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.044317926088339
$ws.Range("D19").Value = 1.052248102849461
$ws.Range("E19").Value = 1.057965713267155
$ws.Range("F19").Value = 1.065135765470662
$ws.Range("I19").Value = 1.048276341941272
$ws.Range("J19").Value = 1.050339392342016
$ws.Range("K19").Value = 1.05550305307239
$ws.Range("L19").Value = 1.061201860645509
$ws.Range("M19").Value = 1.068348646026915
$ws.Range("N19").Value = 1.020637067257782

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.044031870014104
$ws.Range("D20").Value = 1.052024646722475
$ws.Range("E20").Value = 1.057696968931063
$ws.Range("F20").Value = 1.064865525671713
$ws.Range("I20").Value = 1.048196194321076
$ws.Range("J20").Value = 1.050174379937619
$ws.Range("K20").Value = 1.055343142695136
$ws.Range("L20").Value = 1.060996419505753
$ws.Range("M20").Value = 1.068141225475371
$ws.Range("N20").Value = 1.02058179211641

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.043102895841588
$ws.Range("D21").Value = 1.051299022877973
$ws.Range("E21").Value = 1.056824704102395
$ws.Range("F21").Value = 1.063988334188299
$ws.Range("I21").Value = 1.047934349801568
$ws.Range("J21").Value = 1.049637959188787
$ws.Range("K21").Value = 1.054823139365693
$ws.Range("L21").Value = 1.060329026761357
$ws.Range("M21").Value = 1.06746733294682
$ws.Range("N21").Value = 1.020402061051596

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.042519605396804
$ws.Range("D22").Value = 1.050843458502424
$ws.Range("E22").Value = 1.056277397818051
$ws.Range("F22").Value = 1.063437880908181
$ws.Range("I22").Value = 1.047768743504171
$ws.Range("J22").Value = 1.049300735887278
$ws.Range("K22").Value = 1.054496108441158
$ws.Range("L22").Value = 1.059909814667065
$ws.Range("M22").Value = 1.067043984976903
$ws.Range("N22").Value = 1.020289039447359

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.042828742977596
$ws.Range("D23").Value = 1.051084898209749
$ws.Range("E23").Value = 1.05656742827084
$ws.Range("F23").Value = 1.063729584312019
$ws.Range("I23").Value = 1.047856626594483
$ws.Range("J23").Value = 1.049479499554854
$ws.Range("K23").Value = 1.054669481236305
$ws.Range("L23").Value = 1.060132008065657
$ws.Range("M23").Value = 1.067268375528793
$ws.Range("N23").Value = 1.020348955836125

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.044047457690793
$ws.Range("D24").Value = 1.052036823020196
$ws.Range("E24").Value = 1.057711611422433
$ws.Range("F24").Value = 1.0648802499182
$ws.Range("I24").Value = 1.048200567599513
$ws.Range("J24").Value = 1.050183373777563
$ws.Range("K24").Value = 1.055351859089521
$ws.Range("L24").Value = 1.061007615164219
$ws.Range("M24").Value = 1.068152529265811
$ws.Range("N24").Value = 1.020584804997354

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.045465098589518
$ws.Range("D25").Value = 1.05314430664769
$ws.Range("E25").Value = 1.059044185650838
$ws.Range("F25").Value = 1.066220130036718
$ws.Range("I25").Value = 1.048595420679104
$ws.Range("J25").Value = 1.051000334216213
$ws.Range("K25").Value = 1.056143307100973
$ws.Range("L25").Value = 1.062025413220833
$ws.Range("M25").Value = 1.069180029306876
$ws.Range("N25").Value = 1.020858401723068
